$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 12 with the new product data (replacing "Zapatera" entry)
$ws.Range("A12").Value = 24
$ws.Range("B12").Value = "Velador mediano de roble"
$ws.Range("C12").Value = 15
$ws.Range("D12").Value = 140
$ws.Range("E12").Value = "unidad"
$ws.Range("F12").Value = "Dormitorio"

# Delete rows 13 through 22 (old trailing products no longer present)
$ws.Range("A13:F22").EntireRow.Delete()
